# Auto-generated Excel COM-interop script
# Applies numeric cell updates to match target diff across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 8321.727999999999
$ws.Range("I9").Value = 10153.333
$ws.Range("K9").Value = 10153.333
$ws.Range("M9").Value = -9984.333000000001

$ws.Range("H13").Value = 4324.222
$ws.Range("I13").Value = 1177.5
$ws.Range("K13").Value = 1177.5
$ws.Range("M13").Value = -1008.5

$ws.Range("H16").Value = 4999
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H19").Value = 1011.5789
$ws.Range("I19").Value = 571.25
$ws.Range("K19").Value = 571.25
$ws.Range("M19").Value = -396.25

$ws.Range("H26").Value = 2377.5
$ws.Range("I26").Value = 2255
$ws.Range("K26").Value = 2255
$ws.Range("M26").Value = -1911

$ws.Range("H29").Value = 407
$ws.Range("I29").Value = 506.25
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = 1518.75
$ws.Range("L29").Value = 30
$ws.Range("M29").Value = -1237.75
$ws.Range("N29").Value = -592

$ws.Range("H43").Value = 17726.666
$ws.Range("I43").Value = 25559.75
$ws.Range("K43").Value = 25559.75
$ws.Range("M43").Value = -25490.75

$ws.Range("H64").Value = 10156.363
$ws.Range("I64").Value = 3466.3333
$ws.Range("K64").Value = 3466.3333
$ws.Range("M64").Value = -3218.3333

$ws.Range("H67").Value = 10156.363
$ws.Range("I67").Value = 3466.3333
$ws.Range("K67").Value = 3466.3333
$ws.Range("M67").Value = -2608.3333

$ws.Range("H132").Value = 8548642
$ws.Range("I132").Value = 10417854
$ws.Range("J132").Value = 3671.2856
$ws.Range("K132").Value = 31253562
$ws.Range("L132").Value = 11013.8568
$ws.Range("M132").Value = -31251032
$ws.Range("N132").Value = -16073.8568

$ws.Range("H138").Value = 25777.533
$ws.Range("I138").Value = 43300.082
$ws.Range("J138").Value = 5751.7617
$ws.Range("K138").Value = 129900.246
$ws.Range("L138").Value = 17255.2851
$ws.Range("M138").Value = -124760.246
$ws.Range("N138").Value = -27535.2851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 11333.333
$ws.Range("J15").Value = 11333.333
$ws.Range("L15").Value = 11333.333
$ws.Range("N15").Value = -12033.333

$ws.Range("H23").Value = 57500
$ws.Range("J23").Value = 73333.336
$ws.Range("L23").Value = 73333.336
$ws.Range("N23").Value = -73851.336

$ws.Range("H45").Value = 87250.42999999999
$ws.Range("I45").Value = 108097.09
$ws.Range("K45").Value = 108097.09
$ws.Range("M45").Value = -107720.09

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H132").Value = 2598.8394
$ws.Range("I132").Value = 2337.587
$ws.Range("K132").Value = 7012.761
$ws.Range("M132").Value = -4482.761

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 833.3333
$ws.Range("I23").Value = 750
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 750
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = -467
$ws.Range("N23").Value = -1566

$ws.Range("H134").Value = 1561.5306
$ws.Range("I134").Value = 1542.5745
$ws.Range("K134").Value = 4627.7235
$ws.Range("M134").Value = -2092.7235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 4696.231
$ws.Range("I15").Value = 1449
$ws.Range("J15").Value = 5286.636
$ws.Range("K15").Value = 1449
$ws.Range("L15").Value = 5286.636
$ws.Range("M15").Value = -1279
$ws.Range("N15").Value = -5626.636

$ws.Range("H16").Value = 4749.5
$ws.Range("J16").Value = 4749.5
$ws.Range("L16").Value = 4749.5
$ws.Range("N16").Value = -5323.5

$ws.Range("H22").Value = 574.8
$ws.Range("I22").Value = 574.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 574.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -224.8
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 22467.4
$ws.Range("I31").Value = 34102.355
$ws.Range("J31").Value = 3484.0527
$ws.Range("K31").Value = 34102.355
$ws.Range("L31").Value = 3484.0527
$ws.Range("M31").Value = -33807.355
$ws.Range("N31").Value = -4074.0527

$ws.Range("H34").Value = 22467.4
$ws.Range("I34").Value = 34102.355
$ws.Range("J34").Value = 3484.0527
$ws.Range("K34").Value = 34102.355
$ws.Range("L34").Value = 3484.0527
$ws.Range("M34").Value = -33900.355
$ws.Range("N34").Value = -3888.0527

$ws.Range("H37").Value = 20514.25
$ws.Range("J37").Value = 24019
$ws.Range("L37").Value = 24019
$ws.Range("N37").Value = -24233

$ws.Range("H86").Value = 4285.143
$ws.Range("I86").Value = 5398
$ws.Range("K86").Value = 5398
$ws.Range("M86").Value = -4275

$ws.Range("H89").Value = 4285.143
$ws.Range("I89").Value = 5398
$ws.Range("K89").Value = 26990
$ws.Range("M89").Value = -21374

$ws.Range("H113").Value = 4749.5
$ws.Range("J113").Value = 4749.5
$ws.Range("L113").Value = 4749.5
$ws.Range("N113").Value = -9089.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 11532.5
$ws.Range("J106").Value = 12177.111
$ws.Range("L106").Value = 36531.333
$ws.Range("N106").Value = -38423.333

$ws.Range("H109").Value = 556.3333
$ws.Range("I109").Value = 438.375
$ws.Range("J109").Value = 1500
$ws.Range("K109").Value = 1315.125
$ws.Range("L109").Value = 4500
$ws.Range("M109").Value = -275.125
$ws.Range("N109").Value = -6580

$ws.Range("H119").Value = 308
$ws.Range("I119").Value = 308
$ws.Range("K119").Value = 924
$ws.Range("M119").Value = 3914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 29949.5
$ws.Range("J52").Value = 29949.5
$ws.Range("L52").Value = 29949.5
$ws.Range("N52").Value = -30467.5

$ws.Range("H126").Value = 4462.2593
$ws.Range("I126").Value = 3565.8
$ws.Range("J126").Value = 5582.8335
$ws.Range("K126").Value = 10697.4
$ws.Range("L126").Value = 16748.5005
$ws.Range("M126").Value = -8227.400000000001
$ws.Range("N126").Value = -21688.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 3475
$ws.Range("I11").Value = 2900
$ws.Range("J11").Value = 4050
$ws.Range("K11").Value = 2900
$ws.Range("L11").Value = 4050
$ws.Range("M11").Value = -2760
$ws.Range("N11").Value = -4330

$ws.Range("H14").Value = 13000
$ws.Range("I14").Value = 18000
$ws.Range("J14").Value = 8000
$ws.Range("K14").Value = 18000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = -17828
$ws.Range("N14").Value = -8344

$ws.Range("H16").Value = 22728536
$ws.Range("I16").Value = 33334282
$ws.Range("K16").Value = 33334282
$ws.Range("M16").Value = -33334112

$ws.Range("H17").Value = 1832.375
$ws.Range("I17").Value = 1375
$ws.Range("K17").Value = 1375
$ws.Range("M17").Value = -1205

$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040

$ws.Range("H132").Value = 2889.5
$ws.Range("I132").Value = 2889.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8668.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6138.5
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 47850.363
$ws.Range("I136").Value = 60935.766
$ws.Range("J136").Value = 3360
$ws.Range("K136").Value = 182807.298
$ws.Range("L136").Value = 10080
$ws.Range("M136").Value = -180257.298
$ws.Range("N136").Value = -15180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 1672483
$ws.Range("I20").Value = 5000005
$ws.Range("J20").Value = 8722
$ws.Range("K20").Value = 5000005
$ws.Range("L20").Value = 8722
$ws.Range("M20").Value = -4999765
$ws.Range("N20").Value = -9202

$ws.Range("H75").Value = 48000
$ws.Range("I75").Value = 48000
$ws.Range("K75").Value = 48000
$ws.Range("M75").Value = -47064

$ws.Range("H78").Value = 48000
$ws.Range("I78").Value = 48000
$ws.Range("K78").Value = 144000
$ws.Range("M78").Value = -139320

$ws.Range("H122").Value = 3356.8
$ws.Range("I122").Value = 3031.4583
$ws.Range("K122").Value = 9094.374899999999
$ws.Range("M122").Value = -6644.374899999999

$ws.Range("H132").Value = 1199.5883
$ws.Range("I132").Value = 1152.8667
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 3458.6001
$ws.Range("L132").Value = 4650
$ws.Range("M132").Value = -928.6001000000001
$ws.Range("N132").Value = -9710

$ws.Range("H136").Value = 2199.5833
$ws.Range("I136").Value = 2076.2
$ws.Range("J136").Value = 2480
$ws.Range("K136").Value = 6228.599999999999
$ws.Range("L136").Value = 7440
$ws.Range("M136").Value = -3678.599999999999
$ws.Range("N136").Value = -12540
